$wb = $excel.ActiveWorkbook

# Sheet: ARM (29 cell updates)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1495.7
$ws.Range("I2").Value = 1727.75
$ws.Range("K2").Value = 1727.75
$ws.Range("M2").Value = -1614.75
$ws.Range("H32").Value = 3851507.2
$ws.Range("I32").Value = 1855655.8
$ws.Range("J32").Value = 23810022
$ws.Range("K32").Value = 1855655.8
$ws.Range("L32").Value = 23810022
$ws.Range("M32").Value = -1855368.8
$ws.Range("N32").Value = -23810596
$ws.Range("H116").Value = 1495.7
$ws.Range("I116").Value = 1727.75
$ws.Range("K116").Value = 1727.75
$ws.Range("M116").Value = 566.25
$ws.Range("H122").Value = 3206.8948
$ws.Range("I122").Value = 3397.6667
$ws.Range("J122").Value = 2491.5
$ws.Range("K122").Value = 10193.0001
$ws.Range("L122").Value = 7474.5
$ws.Range("M122").Value = -7743.000100000001
$ws.Range("N122").Value = -12374.5
$ws.Range("H132").Value = 2205
$ws.Range("I132").Value = 1076.6666
$ws.Range("J132").Value = 3333.3333
$ws.Range("K132").Value = 3229.9998
$ws.Range("L132").Value = 9999.999899999999
$ws.Range("M132").Value = -699.9998000000001
$ws.Range("N132").Value = -15059.9999

# Sheet: BSM (26 cell updates)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1495.7
$ws.Range("I3").Value = 1727.75
$ws.Range("K3").Value = 1727.75
$ws.Range("M3").Value = -1613.75
$ws.Range("H81").Value = 25884.715
$ws.Range("J81").Value = 27580.666
$ws.Range("L81").Value = 27580.666
$ws.Range("N81").Value = -29702.666
$ws.Range("H84").Value = 25884.715
$ws.Range("J84").Value = 27580.666
$ws.Range("L84").Value = 82741.99800000001
$ws.Range("N84").Value = -93349.99800000001
$ws.Range("H86").Value = 3886.7856
$ws.Range("I86").Value = 3729.2727
$ws.Range("J86").Value = 4464.3335
$ws.Range("K86").Value = 3729.2727
$ws.Range("L86").Value = 4464.3335
$ws.Range("M86").Value = -2606.2727
$ws.Range("N86").Value = -6710.3335
$ws.Range("H89").Value = 3886.7856
$ws.Range("I89").Value = 3729.2727
$ws.Range("J89").Value = 4464.3335
$ws.Range("K89").Value = 18646.3635
$ws.Range("L89").Value = 22321.6675
$ws.Range("M89").Value = -13030.3635
$ws.Range("N89").Value = -33553.6675

# Sheet: CRP (15 cell updates)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 4354.1665
$ws.Range("J99").Value = 4802
$ws.Range("L99").Value = 4802
$ws.Range("N99").Value = -7798
$ws.Range("H122").Value = 3012.875
$ws.Range("I122").Value = 2761
$ws.Range("J122").Value = 4104.3335
$ws.Range("K122").Value = 8283
$ws.Range("L122").Value = 12313.0005
$ws.Range("M122").Value = -5833
$ws.Range("N122").Value = -17213.0005
$ws.Range("H126").Value = 4354.1665
$ws.Range("J126").Value = 4802
$ws.Range("L126").Value = 14406
$ws.Range("N126").Value = -19346

# Sheet: CUL (15 cell updates)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 1619.875
$ws.Range("I114").Value = 1065.5714
$ws.Range("J114").Value = 5500
$ws.Range("K114").Value = 3196.7142
$ws.Range("L114").Value = 16500
$ws.Range("M114").Value = 57.28579999999965
$ws.Range("N114").Value = -23008
$ws.Range("H129").Value = 93854
$ws.Range("J129").Value = 140249.75
$ws.Range("L129").Value = 420749.25
$ws.Range("N129").Value = -430749.25
$ws.Range("H131").Value = 1554925.6
$ws.Range("J131").Value = 1784458.9
$ws.Range("L131").Value = 5353376.699999999
$ws.Range("N131").Value = -5363456.699999999

# Sheet: GSM (22 cell updates)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2032.8823
$ws.Range("I113").Value = 1986
$ws.Range("J113").Value = 2099.8572
$ws.Range("K113").Value = 1986
$ws.Range("L113").Value = 2099.8572
$ws.Range("M113").Value = 184
$ws.Range("N113").Value = -6439.8572
$ws.Range("H126").Value = 13960.182
$ws.Range("I126").Value = 2500
$ws.Range("J126").Value = 15106.2
$ws.Range("K126").Value = 7500
$ws.Range("L126").Value = 45318.60000000001
$ws.Range("M126").Value = -5030
$ws.Range("N126").Value = -50258.60000000001
$ws.Range("H132").Value = 2254
$ws.Range("I132").Value = 2004.8
$ws.Range("K132").Value = 6014.4
$ws.Range("M132").Value = -3484.4
$ws.Range("H134").Value = 49999.5
$ws.Range("J134").Value = 49999.5
$ws.Range("L134").Value = 149998.5
$ws.Range("N134").Value = -155068.5

# Sheet: LTW (47 cell updates)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3140.8572
$ws.Range("I7").Value = 2990.5
$ws.Range("J7").Value = 3201
$ws.Range("K7").Value = 2990.5
$ws.Range("L7").Value = 3201
$ws.Range("M7").Value = -2878.5
$ws.Range("N7").Value = -3425
$ws.Range("H22").Value = 1309.1428
$ws.Range("I22").Value = 1360.6666
$ws.Range("J22").Value = 1000
$ws.Range("K22").Value = 1360.6666
$ws.Range("L22").Value = 1000
$ws.Range("M22").Value = -1065.6666
$ws.Range("N22").Value = -1590
$ws.Range("H27").Value = 1309.1428
$ws.Range("I27").Value = 1360.6666
$ws.Range("J27").Value = 1000
$ws.Range("K27").Value = 1360.6666
$ws.Range("L27").Value = 1000
$ws.Range("M27").Value = -1253.6666
$ws.Range("N27").Value = -1214
$ws.Range("H40").Value = 24362.545
$ws.Range("I40").Value = 24362.545
$ws.Range("K40").Value = 24362.545
$ws.Range("M40").Value = -24226.545
$ws.Range("H61").Value = 1932.6666
$ws.Range("I61").Value = 1932.6666
$ws.Range("K61").Value = 1932.6666
$ws.Range("M61").Value = -1730.6666
$ws.Range("H113").Value = 1932.6666
$ws.Range("I113").Value = 1932.6666
$ws.Range("K113").Value = 1932.6666
$ws.Range("M113").Value = 237.3334
$ws.Range("H122").Value = 4887.923
$ws.Range("I122").Value = 2093.8572
$ws.Range("J122").Value = 8147.6665
$ws.Range("K122").Value = 6281.571599999999
$ws.Range("L122").Value = 24442.9995
$ws.Range("M122").Value = -3831.571599999999
$ws.Range("N122").Value = -29342.9995
$ws.Range("H126").Value = 3140.8572
$ws.Range("I126").Value = 2990.5
$ws.Range("J126").Value = 3201
$ws.Range("K126").Value = 8971.5
$ws.Range("L126").Value = 9603
$ws.Range("M126").Value = -6501.5
$ws.Range("N126").Value = -14543

# Sheet: WVR (7 cell updates)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 699.9375
$ws.Range("J107").Value = 724.6667
$ws.Range("L107").Value = 2174.0001
$ws.Range("N107").Value = -6014.0001
$ws.Range("I122").Value = 2032.5333
$ws.Range("K122").Value = 6097.5999
$ws.Range("M122").Value = -3647.5999
